# Generate Report for Handback
# Adds a second handed-back file ("f4c4224a-...") to the Overview / zh-cn / de-de
# sheets, and refreshes the existing row's UUID + timestamps for the first file
# (renamed from "2cdd9b1e-..." to "d6daab99-...").

$wb = $excel.ActiveWorkbook

$oldGuid = "2cdd9b1e-293e-4ac6-a36a-b2e9d06c37a3"
$newGuid1 = "d6daab99-a7b2-49d3-934c-621dbeee224a"
$newGuid2 = "f4c4224a-7307-4a7e-88f5-81b7a6155c24"

$oldZhHash = "04104383c7affb9ef1f9a4f05e6882cb10d5b276"
$newHash1 = "42a791ae6569292e1e90ab79ea360b97de3b587a"
$newHash2 = "b024e102f4bc90e294475d4cfcca555621c4e720"

$srcRepoUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/fdcd04a489f5c3620c8420509f597f2c31982caa/e2e/"
$zhRepoUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/48e9fb91bcf87a3d17eb0c34ee0024e97944d0da/e2e/"
$deRepoUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/1d0dd8235f8f551fddfac2a82888ec02b4a4c0cb/e2e/"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsO = $wb.Worksheets.Item("Overview")

# refresh row 2 (existing file renamed + new generate date)
$wsO.Range("A2").Value = "$newGuid1.md"
$wsO.Range("B2").Hyperlinks.Delete()
$wsO.Hyperlinks.Add($wsO.Range("B2"), "$srcRepoUrl$newGuid1.md", "", "", "e2e\$newGuid1.md")
$wsO.Range("G2").Value = "2016-08-13 09:14:34"

# add row 3 for the new file
$loO = $wsO.ListObjects.Item(1)
$loO.ListRows.Add() | Out-Null
$wsO.Range("A3").Value = "$newGuid2.md"
$wsO.Range("B3").Value = "e2e\$newGuid2.md"
$wsO.Hyperlinks.Add($wsO.Range("B3"), "$srcRepoUrl$newGuid2.md", "", "", "e2e\$newGuid2.md")
$wsO.Range("C3").Value = ".md"
$wsO.Range("E3").Value = "Handed back: in sync with en-US"
$wsO.Range("F3").Value = "Handed back: in sync with en-US"
$wsO.Range("G3").Value = "2016-08-13 09:14:34"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("A2").Value = "$newGuid1.md"
$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "$srcRepoUrl$newGuid1.md", "", "", "$newGuid1.md")
$wsZh.Range("G2").Value = "$newGuid1.$newHash1.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-13 09:14:27"
$wsZh.Range("I2").Value = "$newGuid1.md"
$wsZh.Range("I2").Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "$zhRepoUrl$newGuid1.md", "", "", "$newGuid1.md")
$wsZh.Range("J2").Value = "$newGuid1.$newHash1.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-13 09:14:55"

$loZh = $wsZh.ListObjects.Item(1)
$loZh.ListRows.Add() | Out-Null
$wsZh.Range("A3").Value = "$newGuid2.md"
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "$srcRepoUrl$newGuid2.md", "", "", "$newGuid2.md")
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = "$newGuid2.$newHash2.zh-cn.xlf"
$wsZh.Range("H3").Value = "2016-08-13 09:14:27"
$wsZh.Range("I3").Value = "$newGuid2.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "$zhRepoUrl$newGuid2.md", "", "", "$newGuid2.md")
$wsZh.Range("J3").Value = "$newGuid2.$newHash2.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-13 09:14:55"
$wsZh.Range("L3").Value = ""
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N3").Value = ""
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P3").Value = ""

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("A2").Value = "$newGuid1.md"
$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "$srcRepoUrl$newGuid1.md", "", "", "$newGuid1.md")
$wsDe.Range("G2").Value = "$newGuid1.$newHash1.de-de.xlf"
$wsDe.Range("H2").Value = "2016-08-13 09:14:34"
$wsDe.Range("I2").Value = "$newGuid1.md"
$wsDe.Range("I2").Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "$deRepoUrl$newGuid1.md", "", "", "$newGuid1.md")
$wsDe.Range("J2").Value = "$newGuid1.$newHash1.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-13 09:15:10"

$loDe = $wsDe.ListObjects.Item(1)
$loDe.ListRows.Add() | Out-Null
$wsDe.Range("A3").Value = "$newGuid2.md"
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "$srcRepoUrl$newGuid2.md", "", "", "$newGuid2.md")
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = "$newGuid2.$newHash2.de-de.xlf"
$wsDe.Range("H3").Value = "2016-08-13 09:14:34"
$wsDe.Range("I3").Value = "$newGuid2.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "$deRepoUrl$newGuid2.md", "", "", "$newGuid2.md")
$wsDe.Range("J3").Value = "$newGuid2.$newHash2.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-13 09:15:10"
$wsDe.Range("L3").Value = ""
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N3").Value = ""
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P3").Value = ""
